# New weekly price record for Puerro (Vega Modelo de Temuco) is inserted as
# row 57, pushing the existing rows 57-137 down to 58-138 (dimension grows
# from A1:R137 to A1:R138).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(57).EntireRow.Insert()

$ws.Cells.Item(57, 1).Value = 10
$ws.Cells.Item(57, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(57, 3).Value = 'La Araucanía'
$ws.Cells.Item(57, 4).Value = 44482
$ws.Cells.Item(57, 5).Value = 9
$ws.Cells.Item(57, 6).Value = 100112005
$ws.Cells.Item(57, 7).Value = 'Puerro'
$ws.Cells.Item(57, 8).Value = 'Azul de Maquehue'
$ws.Cells.Item(57, 9).Value = 'Primera'
$ws.Cells.Item(57, 10).Value = 40
$ws.Cells.Item(57, 11).Value = 7000
$ws.Cells.Item(57, 12).Value = 7000
$ws.Cells.Item(57, 13).Value = 7000
$ws.Cells.Item(57, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(57, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(57, 16).Value = 583
$ws.Cells.Item(57, 17).Value = 12
$ws.Cells.Item(57, 18).Value = 'Hortaliza'
